$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
  2  = @{ D = "307.69";     E = "-6.06%" }
  3  = @{ D = "40.69";      E = "-8.14%" }
  4  = @{ D = "5.055";      E = "-4.57%" }
  5  = @{ D = "0.07790";    E = "-6.86%" }
  6  = @{ D = "4.318";      E = "-1.95%" }
  7  = @{ D = "1.650";      E = "-15.46%" }
  8  = @{ D = "0.9075";     E = "-6.45%" }
  9  = @{ D = "0.1050";     E = "-7.46%" }
  10 = @{ D = "0.1746";     E = "-8.16%" }
  11 = @{ D = "0.04474";    E = "-2.80%" }
  12 = @{ D = "0.08967";    E = "-7.10%" }
  13 = @{ E = "-15.62%" }
  14 = @{ D = "0.1056";     E = "-0.36%" }
  15 = @{ D = "0.001258";   E = "-3.19%" }
  16 = @{ D = "0.005719";   E = "-1.50%" }
  17 = @{ D = "3.370";      E = "-0.54%" }
  18 = @{ D = "2.559";      E = "1.94%" }
  19 = @{ D = "0.3366";     E = "0.23%" }
  20 = @{ E = "-0.18%" }
  21 = @{ D = "0.2858";     E = "10.94%" }
  22 = @{ D = "0.04183";    E = "0.78%" }
  23 = @{ E = "-0.93%" }
  24 = @{ D = "0.004114";   E = "-6.89%" }
  25 = @{ D = "0.0001232";  E = "-5.30%" }
  26 = @{ D = "0.0003000";  E = "0.61%" }
  38 = @{ D = "0.02412";    E = "-11.12%" }
  39 = @{ D = "0.05183";    E = "-7.62%" }
  40 = @{ D = "0.008026";   E = "2.13%" }
  41 = @{ D = "0.1328";     E = "-5.97%" }
  42 = @{ D = "0.007514";   E = "2.70%" }
  43 = @{ D = "0.001998";   E = "-2.57%" }
  44 = @{ D = "0.008090";   E = "-6.59%" }
  45 = @{ D = "0.3345";     E = "-4.86%" }
  46 = @{ D = "0.00006738"; E = "-1.52%" }
  47 = @{ D = "0.00000000755"; E = "0.68%" }
  48 = @{ D = "0.003361";   E = "-3.79%" }
  49 = @{ D = "0.004129";   E = "16.86%" }
  50 = @{ D = "0.00002115"; E = "0.68%" }
  51 = @{ D = "0.0002015";  E = "0.68%" }
}

foreach ($row in $updates.Keys) {
  $cols = $updates[$row]
  foreach ($col in $cols.Keys) {
    $ws.Range("$col$row").Value = "'" + $cols[$col]
  }
}
